$d = $word.ActiveDocument

$d.Content.Find.Execute("11+45=", $true, $false, $false, $false, $false, $true, 1, $false, "54-15=", 2) | Out-Null
$d.Content.Find.Execute("35+37=", $true, $false, $false, $false, $false, $true, 1, $false, "38+58=", 2) | Out-Null
$d.Content.Find.Execute("20+41=", $true, $false, $false, $false, $false, $true, 1, $false, "48-1=", 2) | Out-Null
$d.Content.Find.Execute("92-3=", $true, $false, $false, $false, $false, $true, 1, $false, "36-26=", 2) | Out-Null
$d.Content.Find.Execute("28-14=", $true, $false, $false, $false, $false, $true, 1, $false, "84-6=", 2) | Out-Null
$d.Content.Find.Execute("73+19=", $true, $false, $false, $false, $false, $true, 1, $false, "98-83=", 2) | Out-Null
$d.Content.Find.Execute("42+3=", $true, $false, $false, $false, $false, $true, 1, $false, "73-60=", 2) | Out-Null
$d.Content.Find.Execute("17+44=", $true, $false, $false, $false, $false, $true, 1, $false, "46+13=", 2) | Out-Null
$d.Content.Find.Execute("48-23=", $true, $false, $false, $false, $false, $true, 1, $false, "90-28=", 2) | Out-Null
$d.Content.Find.Execute("63+30=", $true, $false, $false, $false, $false, $true, 1, $false, "72-67=", 2) | Out-Null
$d.Content.Find.Execute("74-35=", $true, $false, $false, $false, $false, $true, 1, $false, "29+59=", 2) | Out-Null
$d.Content.Find.Execute("69-11=", $true, $false, $false, $false, $false, $true, 1, $false, "31+43=", 2) | Out-Null
$d.Content.Find.Execute("15+6=", $true, $false, $false, $false, $false, $true, 1, $false, "0+61=", 2) | Out-Null
$d.Content.Find.Execute("81-9=", $true, $false, $false, $false, $false, $true, 1, $false, "21+28=", 2) | Out-Null
$d.Content.Find.Execute("49-9=", $true, $false, $false, $false, $false, $true, 1, $false, "62-43=", 2) | Out-Null
$d.Content.Find.Execute("1+52=", $true, $false, $false, $false, $false, $true, 1, $false, "44+46=", 2) | Out-Null
$d.Content.Find.Execute("25-1=", $true, $false, $false, $false, $false, $true, 1, $false, "74-57=", 2) | Out-Null
$d.Content.Find.Execute("84-75=", $true, $false, $false, $false, $false, $true, 1, $false, "2+26=", 2) | Out-Null
$d.Content.Find.Execute("34+49=", $true, $false, $false, $false, $false, $true, 1, $false, "84-61=", 2) | Out-Null
$d.Content.Find.Execute("89-28=", $true, $false, $false, $false, $false, $true, 1, $false, "32-11=", 2) | Out-Null
$d.Content.Find.Execute("87-78=", $true, $false, $false, $false, $false, $true, 1, $false, "3+87=", 2) | Out-Null
$d.Content.Find.Execute("38+40=", $true, $false, $false, $false, $false, $true, 1, $false, "12+8=", 2) | Out-Null
$d.Content.Find.Execute("69-52=", $true, $false, $false, $false, $false, $true, 1, $false, "71-10=", 2) | Out-Null
$d.Content.Find.Execute("35+57=", $true, $false, $false, $false, $false, $true, 1, $false, "41+56=", 2) | Out-Null
$d.Content.Find.Execute("60-29=", $true, $false, $false, $false, $false, $true, 1, $false, "37-30=", 2) | Out-Null
$d.Content.Find.Execute("22-14=", $true, $false, $false, $false, $false, $true, 1, $false, "39+21=", 2) | Out-Null
$d.Content.Find.Execute("65+19=", $true, $false, $false, $false, $false, $true, 1, $false, "33-10=", 2) | Out-Null
$d.Content.Find.Execute("43-0=", $true, $false, $false, $false, $false, $true, 1, $false, "48+50=", 2) | Out-Null
$d.Content.Find.Execute("32+27=", $true, $false, $false, $false, $false, $true, 1, $false, "4+31=", 2) | Out-Null
$d.Content.Find.Execute("78-70=", $true, $false, $false, $false, $false, $true, 1, $false, "97-75=", 2) | Out-Null
$d.Content.Find.Execute("71-68=", $true, $false, $false, $false, $false, $true, 1, $false, "85-10=", 2) | Out-Null
$d.Content.Find.Execute("36+46=", $true, $false, $false, $false, $false, $true, 1, $false, "51-9=", 2) | Out-Null
$d.Content.Find.Execute("70-18=", $true, $false, $false, $false, $false, $true, 1, $false, "71-53=", 2) | Out-Null
$d.Content.Find.Execute("37+30=", $true, $false, $false, $false, $false, $true, 1, $false, "61-16=", 2) | Out-Null
$d.Content.Find.Execute("86-5=", $true, $false, $false, $false, $false, $true, 1, $false, "77-50=", 2) | Out-Null
$d.Content.Find.Execute("42-32=", $true, $false, $false, $false, $false, $true, 1, $false, "7+35=", 2) | Out-Null
$d.Content.Find.Execute("75+5=", $true, $false, $false, $false, $false, $true, 1, $false, "72+20=", 2) | Out-Null
$d.Content.Find.Execute("98-39=", $true, $false, $false, $false, $false, $true, 1, $false, "53-25=", 2) | Out-Null
$d.Content.Find.Execute("60+33=", $true, $false, $false, $false, $false, $true, 1, $false, "23+58=", 2) | Out-Null
$d.Content.Find.Execute("35+58=", $true, $false, $false, $false, $false, $true, 1, $false, "25-18=", 2) | Out-Null
$d.Content.Find.Execute("71+12=", $true, $false, $false, $false, $false, $true, 1, $false, "88-5=", 2) | Out-Null
$d.Content.Find.Execute("80-78=", $true, $false, $false, $false, $false, $true, 1, $false, "61+35=", 2) | Out-Null
$d.Content.Find.Execute("85+4=", $true, $false, $false, $false, $false, $true, 1, $false, "63+36=", 2) | Out-Null
$d.Content.Find.Execute("46+8=", $true, $false, $false, $false, $false, $true, 1, $false, "47-8=", 2) | Out-Null
$d.Content.Find.Execute("68-16=", $true, $false, $false, $false, $false, $true, 1, $false, "58+29=", 2) | Out-Null
$d.Content.Find.Execute("49+28=", $true, $false, $false, $false, $false, $true, 1, $false, "88-87=", 2) | Out-Null
$d.Content.Find.Execute("51+43=", $true, $false, $false, $false, $false, $true, 1, $false, "56-4=", 2) | Out-Null
$d.Content.Find.Execute("28+33=", $true, $false, $false, $false, $false, $true, 1, $false, "76-32=", 2) | Out-Null
$d.Content.Find.Execute("61-61=", $true, $false, $false, $false, $false, $true, 1, $false, "47+50=", 2) | Out-Null
$d.Content.Find.Execute("87-68=", $true, $false, $false, $false, $false, $true, 1, $false, "48+6=", 2) | Out-Null
$d.Content.Find.Execute("50+23=", $true, $false, $false, $false, $false, $true, 1, $false, "33+34=", 2) | Out-Null
$d.Content.Find.Execute("72+17=", $true, $false, $false, $false, $false, $true, 1, $false, "61-4=", 2) | Out-Null
$d.Content.Find.Execute("39-20=", $true, $false, $false, $false, $false, $true, 1, $false, "81-10=", 2) | Out-Null
$d.Content.Find.Execute("40-23=", $true, $false, $false, $false, $false, $true, 1, $false, "1+45=", 2) | Out-Null
$d.Content.Find.Execute("81-78=", $true, $false, $false, $false, $false, $true, 1, $false, "46-2=", 2) | Out-Null
$d.Content.Find.Execute("59-32=", $true, $false, $false, $false, $false, $true, 1, $false, "95-58=", 2) | Out-Null
$d.Content.Find.Execute("56-40=", $true, $false, $false, $false, $false, $true, 1, $false, "41+25=", 2) | Out-Null
$d.Content.Find.Execute("90-79=", $true, $false, $false, $false, $false, $true, 1, $false, "71-10=", 2) | Out-Null
$d.Content.Find.Execute("17+51=", $true, $false, $false, $false, $false, $true, 1, $false, "69-17=", 2) | Out-Null
$d.Content.Find.Execute("50-6=", $true, $false, $false, $false, $false, $true, 1, $false, "23+54=", 2) | Out-Null
$d.Content.Find.Execute("43+35=", $true, $false, $false, $false, $false, $true, 1, $false, "23+57=", 2) | Out-Null
$d.Content.Find.Execute("73-71=", $true, $false, $false, $false, $false, $true, 1, $false, "26+53=", 2) | Out-Null
$d.Content.Find.Execute("45-28=", $true, $false, $false, $false, $false, $true, 1, $false, "51+33=", 2) | Out-Null
$d.Content.Find.Execute("74+4=", $true, $false, $false, $false, $false, $true, 1, $false, "95-39=", 2) | Out-Null
$d.Content.Find.Execute("7+47=", $true, $false, $false, $false, $false, $true, 1, $false, "43-17=", 2) | Out-Null
$d.Content.Find.Execute("48-9=", $true, $false, $false, $false, $false, $true, 1, $false, "52-8=", 2) | Out-Null
$d.Content.Find.Execute("95-80=", $true, $false, $false, $false, $false, $true, 1, $false, "52-39=", 2) | Out-Null
$d.Content.Find.Execute("29+69=", $true, $false, $false, $false, $false, $true, 1, $false, "64+23=", 2) | Out-Null
$d.Content.Find.Execute("9+9=", $true, $false, $false, $false, $false, $true, 1, $false, "10+61=", 2) | Out-Null
$d.Content.Find.Execute("20-7=", $true, $false, $false, $false, $false, $true, 1, $false, "13+61=", 2) | Out-Null
$d.Content.Find.Execute("21-18=", $true, $false, $false, $false, $false, $true, 1, $false, "70+28=", 2) | Out-Null
$d.Content.Find.Execute("39-24=", $true, $false, $false, $false, $false, $true, 1, $false, "42+18=", 2) | Out-Null
$d.Content.Find.Execute("95-63=", $true, $false, $false, $false, $false, $true, 1, $false, "85-83=", 2) | Out-Null
$d.Content.Find.Execute("73-24=", $true, $false, $false, $false, $false, $true, 1, $false, "7-3=", 2) | Out-Null
$d.Content.Find.Execute("55-42=", $true, $false, $false, $false, $false, $true, 1, $false, "47+24=", 2) | Out-Null
$d.Content.Find.Execute("45+26=", $true, $false, $false, $false, $false, $true, 1, $false, "3+42=", 2) | Out-Null
$d.Content.Find.Execute("13+46=", $true, $false, $false, $false, $false, $true, 1, $false, "6+80=", 2) | Out-Null
$d.Content.Find.Execute("71+11=", $true, $false, $false, $false, $false, $true, 1, $false, "30+12=", 2) | Out-Null
$d.Content.Find.Execute("81-16=", $true, $false, $false, $false, $false, $true, 1, $false, "68+15=", 2) | Out-Null
$d.Content.Find.Execute("4+33=", $true, $false, $false, $false, $false, $true, 1, $false, "42+39=", 2) | Out-Null
$d.Content.Find.Execute("9+24=", $true, $false, $false, $false, $false, $true, 1, $false, "99-98=", 2) | Out-Null
$d.Content.Find.Execute("24+42=", $true, $false, $false, $false, $false, $true, 1, $false, "49-7=", 2) | Out-Null
$d.Content.Find.Execute("73-53=", $true, $false, $false, $false, $false, $true, 1, $false, "8+81=", 2) | Out-Null
$d.Content.Find.Execute("19+5=", $true, $false, $false, $false, $false, $true, 1, $false, "83+14=", 2) | Out-Null
$d.Content.Find.Execute("79-36=", $true, $false, $false, $false, $false, $true, 1, $false, "91-50=", 2) | Out-Null
$d.Content.Find.Execute("2+34=", $true, $false, $false, $false, $false, $true, 1, $false, "4+21=", 2) | Out-Null
$d.Content.Find.Execute("86-51=", $true, $false, $false, $false, $false, $true, 1, $false, "60-16=", 2) | Out-Null
$d.Content.Find.Execute("50-35=", $true, $false, $false, $false, $false, $true, 1, $false, "82-16=", 2) | Out-Null
$d.Content.Find.Execute("56-43=", $true, $false, $false, $false, $false, $true, 1, $false, "39-3=", 2) | Out-Null
$d.Content.Find.Execute("74-68=", $true, $false, $false, $false, $false, $true, 1, $false, "56+7=", 2) | Out-Null
$d.Content.Find.Execute("7+57=", $true, $false, $false, $false, $false, $true, 1, $false, "49+30=", 2) | Out-Null
$d.Content.Find.Execute("63-4=", $true, $false, $false, $false, $false, $true, 1, $false, "33+48=", 2) | Out-Null
$d.Content.Find.Execute("93-67=", $true, $false, $false, $false, $false, $true, 1, $false, "92-58=", 2) | Out-Null
$d.Content.Find.Execute("25-21=", $true, $false, $false, $false, $false, $true, 1, $false, "80-48=", 2) | Out-Null
$d.Content.Find.Execute("59+7=", $true, $false, $false, $false, $false, $true, 1, $false, "15+31=", 2) | Out-Null
$d.Content.Find.Execute("34+24=", $true, $false, $false, $false, $false, $true, 1, $false, "64-16=", 2) | Out-Null
$d.Content.Find.Execute("32-7=", $true, $false, $false, $false, $false, $true, 1, $false, "73-33=", 2) | Out-Null
$d.Content.Find.Execute("84-57=", $true, $false, $false, $false, $false, $true, 1, $false, "26+14=", 2) | Out-Null
$d.Content.Find.Execute("50+13=", $true, $false, $false, $false, $false, $true, 1, $false, "16+17=", 2) | Out-Null
$d.Content.Find.Execute("69-50=", $true, $false, $false, $false, $false, $true, 1, $false, "12-1=", 2) | Out-Null
